# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a value into a "Price" (column D) cell as literal text,
# preserving the cell's original style (Excel would otherwise silently
# parse strings like "5.060" or "0.07210" as numbers and drop the formatting
# / trailing zeros that the source data relies on).
function Set-PriceText($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-PriceText $ws.Range("D2") "26.985.68"
$ws.Range("E2").Value = "  +0.08%  "
Set-PriceText $ws.Range("D3") "1.844.64"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.48%  "
Set-PriceText $ws.Range("D5") "1.012"
$ws.Range("E5").Value = "  +0.38%  "
Set-PriceText $ws.Range("D6") "309.14"
$ws.Range("E6").Value = "  -0.33%  "
Set-PriceText $ws.Range("D7") "0.4761"
$ws.Range("E7").Value = "  +2.13%  "
Set-PriceText $ws.Range("D8") "0.3677"
$ws.Range("E8").Value = "  +1.36%  "
Set-PriceText $ws.Range("D9") "0.07210"
$ws.Range("E9").Value = "  +0.53%  "
Set-PriceText $ws.Range("D10") "0.9328"
$ws.Range("E10").Value = "  +1.44%  "
Set-PriceText $ws.Range("D11") "19.90"
$ws.Range("E11").Value = "  +1.62%  "
Set-PriceText $ws.Range("D12") "0.07730"
$ws.Range("E12").Value = "  +0.49%  "
Set-PriceText $ws.Range("D13") "1.863.33"
$ws.Range("E13").Value = "  +0.26%  "
Set-PriceText $ws.Range("D14") "5.385"
$ws.Range("E14").Value = "  +1.66%  "
Set-PriceText $ws.Range("D15") "6.471"
$ws.Range("E15").Value = "  +0.98%  "
Set-PriceText $ws.Range("D16") "88.88"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("E17").Value = "  +0.42%  "
Set-PriceText $ws.Range("D18") "0.000008654"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +0.43%  "
Set-PriceText $ws.Range("D20") "27.013.16"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +1.01%  "
Set-PriceText $ws.Range("D22") "5.060"
$ws.Range("E22").Value = "  +0.60%  "
Set-PriceText $ws.Range("D23") "10.63"
$ws.Range("E23").Value = "  +0.00%  "
Set-PriceText $ws.Range("D24") "1.943"
$ws.Range("E24").Value = "  +0.74%  "
Set-PriceText $ws.Range("D25") "152.67"
$ws.Range("E25").Value = "  +0.39%  "
Set-PriceText $ws.Range("D26") "18.23"
$ws.Range("E26").Value = "  +0.12%  "
Set-PriceText $ws.Range("D27") "2.004"
$ws.Range("E27").Value = "  -3.14%  "
Set-PriceText $ws.Range("D28") "114.34"
$ws.Range("E28").Value = "  +0.26%  "
Set-PriceText $ws.Range("D29") "4.974"
$ws.Range("E29").Value = "  +1.15%  "
Set-PriceText $ws.Range("D30") "0.08863"
$ws.Range("E30").Value = "  +0.09%  "
Set-PriceText $ws.Range("D31") "3.294"
$ws.Range("E31").Value = "  +3.36%  "
Set-PriceText $ws.Range("D33") "0.7397"
$ws.Range("E33").Value = "  -0.54%  "
Set-PriceText $ws.Range("D34") "4.507"
$ws.Range("E34").Value = "  +0.84%  "
Set-PriceText $ws.Range("D35") "2.667"
$ws.Range("E35").Value = "  -6.60%  "
Set-PriceText $ws.Range("D36") "1.112"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("E37").Value = "  +1.26%  "
Set-PriceText $ws.Range("D38") "0.05269"
$ws.Range("E38").Value = "  +1.94%  "
Set-PriceText $ws.Range("D39") "2.972"
$ws.Range("E39").Value = "  +0.31%  "
Set-PriceText $ws.Range("D40") "0.5258"
$ws.Range("E40").Value = "  +2.09%  "
Set-PriceText $ws.Range("D41") "7.031"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("E42").Value = "  +0.19%  "
Set-PriceText $ws.Range("D43") "8.292"
$ws.Range("E43").Value = "  +1.42%  "
Set-PriceText $ws.Range("D44") "10.62"
$ws.Range("E44").Value = "  +0.87%  "
Set-PriceText $ws.Range("D45") "0.4737"
$ws.Range("E45").Value = "  +0.84%  "
Set-PriceText $ws.Range("D46") "1.013"
$ws.Range("E46").Value = "  +0.42%  "
Set-PriceText $ws.Range("D47") "101.75"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("E48").Value = "  +0.57%  "
Set-PriceText $ws.Range("D49") "65.75"
$ws.Range("E49").Value = "  +1.77%  "
Set-PriceText $ws.Range("D50") "0.06064"
$ws.Range("E50").Value = "  +0.36%  "
Set-PriceText $ws.Range("D51") "0.8918"
$ws.Range("E51").Value = "  +3.46%  "
